# Auto-generated edit script: updates cryptos.xlsx cell values
# per the commit diff (price/volume refresh + one two-row coin swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''60.619.68'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -2.64%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.903.42'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -3.75%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.09%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''586.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -1.45%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''147.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -0.77%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.10%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  -2.70%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''2.902.57'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -3.71%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''6.71'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +4.71%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  -4.02%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '''  -2.40%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  -3.73%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''34.06'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -1.17%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = '''  +0.45%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''3.385.04'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -3.70%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''6.83'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -2.55%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''60.595.36'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -2.60%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''2.906.10'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -3.69%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''427.72'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -4.64%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''13.63'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -4.04%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''0.669'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -3.03%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''7.10'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -4.13%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''80.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -2.14%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''11.09'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +1.62%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''2.21'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -1.59%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''11.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -1.86%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  -0.01%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  +0.08%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '''  +0.54%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = '''  -3.04%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  +2.26%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''26.51'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -3.73%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  -3.10%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''0.0₃0839'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -1.45%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  -2.21%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''5.67'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -3.05%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = '''  -1.62%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = '''  -0.66%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''49.27'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -1.80%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = '''Kaspa'
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = '''0.122'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -1.22%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = '''Cosmos'
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = '''8.74'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -3.63%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''0.291'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +1.19%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''41.53'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +1.35%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '''  -1.31%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''370.53'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = '''133.75'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -0.60%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''2.653.03'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -2.98%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D50").Value = '''25.13'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +5.73%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  -1.31%  '
$ws.Range("E51").Style = "Normal"
